$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4264
$ws.Range("C2").Value = 2132

$ws.Range("B3").Value = 73.87
$ws.Range("C3").Value = 36.93

$ws.Range("B4").Value = 160
$ws.Range("C4").Value = 80

$ws.Range("B5").Value = 26.67
$ws.Range("C5").Value = 13.33

$ws.Range("B6").Value = 4003.47
$ws.Range("C6").Value = 2001.73
